$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching the style of the existing
# header row (B1:H1 use style index 1 - bold header style). Copy the format
# from the adjacent "IP" header cell (H1) so the same style is reused rather
# than a new one being created.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data for columns I (I0) and J (IF), one row per data row (rows 2-77).
$data = @(
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 9),
    @(7, 8),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(6, 7),
    @(5, 6),
    @(8, 9),
    @(9, 10),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(5, 6),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(7, 7),
    @(4, 4),
    @(3, 3)
)

$startRow = 2
for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $startRow + $idx
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
